$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.701.14"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").Value = "3.186.46"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "534.67"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").Value = "144.49"
$ws.Range("E6").Value = "  +3.02%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "0.530"
$ws.Range("E8").Value = "  +2.30%  "
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("E10").Value = "  +1.76%  "
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("D12").Value = "3.738.64"
$ws.Range("E12").Value = "  +1.18%  "
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D14").Value = "26.02"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").Value = "59.759.33"
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("D17").Value = "3.188.97"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("D18").Value = "6.21"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "13.12"
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("D20").Value = "8.18"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").Value = "366.57"
$ws.Range("E21").Value = "  -1.58%  "
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("D24").Value = "69.59"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "8.77"
$ws.Range("E25").Value = "  +9.87%  "
$ws.Range("D26").Value = "0.167"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").Value = "0.0₃0882"
$ws.Range("E28").Value = "  +1.71%  "
$ws.Range("D29").Value = "22.27"
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("D32").Value = "5.28"
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("D33").Value = "1.19"
$ws.Range("E33").Value = "  +2.04%  "
$ws.Range("D34").Value = "6.55"
$ws.Range("E34").Value = "  +4.76%  "
$ws.Range("D35").Value = "157.10"
$ws.Range("E35").Value = "  -1.57%  "
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("D37").Value = "2.776.70"
$ws.Range("E37").Value = "  +4.22%  "
$ws.Range("D38").Value = "25.65"
$ws.Range("E38").Value = "  +1.67%  "
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").Value = "0.0291"
$ws.Range("E42").Value = "  +1.35%  "
$ws.Range("D43").Value = "39.32"
$ws.Range("E43").Value = "  +2.09%  "
$ws.Range("D44").Value = "0.712"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("E45").Value = "  +1.19%  "
$ws.Range("D46").Value = "3.228.01"
$ws.Range("E46").Value = "  +1.05%  "
$ws.Range("D47").Value = "0.980"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").Value = "0.810"
$ws.Range("E48").Value = "  +6.56%  "
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("D50").Value = "20.39"
$ws.Range("E50").Value = "  +1.10%  "
$ws.Range("E51").Value = "  +0.00%  "
